$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.598999999999999
$ws.Range("D5").Value = -8.148999999999999
$ws.Range("D9").Value = -7.634
$ws.Range("D11").Value = -8.359
$ws.Range("B21").Value = 6.343
$ws.Range("D21").Value = -7.775999999999999
$ws.Range("B23").Value = 6.842000000000001
$ws.Range("B25").Value = 6.556999999999999
